# Add the newest arabidopsis results for the bnstruct package.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the existing "Global MIT" row (old row 8),
# pushing it down to row 10, and leaving room for two new bnstruct rows.
$ws.Rows.Item(8).Insert(-4121)
$ws.Rows.Item(8).Insert(-4121)

# New row 8: bnstruct, arabidopsis, time step 3
$ws.Range("A8").Value = "bnstruct"
$ws.Range("B8").Value = "arabidopsis"
$ws.Range("C8").Value = 3
$ws.Range("D8").Value = 0.533
$ws.Range("E8").Value = 0.475
$ws.Range("F8").Value = 0.368

# New row 9: bnstruct, arabidopsis, time step 4
$ws.Range("A9").Value = "bnstruct"
$ws.Range("B9").Value = "arabidopsis"
$ws.Range("C9").Value = 4
$ws.Range("D9").Value = 0.491
$ws.Range("E9").Value = 0.346
$ws.Range("F9").Value = 0.451

# New footnote-style summary rows 12-17 (row 11 intentionally left blank)
$ws.Range("A12").Value = "bnstruct, arabidopsis, 3, 0.533, 0.475, 0.368"
$ws.Range("A13").Value = "bnstruct, arabidopsis, 4, 0.491, 0.436, 0.451"
$ws.Range("A14").Value = "bnstruct, arabidopsis, 5, 0.512, 0.458, 0.469"
$ws.Range("A15").Value = "bnstruct, arabidopsis, 6, 0.472, 0.425, 0.497"
$ws.Range("A16").Value = "bnstruct, arabidopsis, 7, 0.573, 0.577, 0.189"
$ws.Range("A17").Value = "bnstruct, arabidopsis, 8, 0.557, 0.532, 0.136"

# Style the footnote rows: 9pt black Times New Roman, left/center aligned, indented.
$footnote = $ws.Range("A12")
$footnote.Font.Size = 9
$footnote.Font.Color = 0
$footnote.HorizontalAlignment = -4131
$footnote.VerticalAlignment = -4108
$footnote.IndentLevel = 1

$footnote.Copy()
$ws.Range("A13:A17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the sheet selection to reflect the new footnote block.
$ws.Range("A12:A17").Select()
